# Calculated_variables.xlsx - "Add files via upload" re-upload of recomputed data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: N7/P7/R7 were placeholder "same" text; replace with the actual
#     numeric Max_Mass_g / Max_SL_mm / Max_TL_mm values (P7 corrected to 144.9).
$ws.Range("N7").Value = 35.24
$ws.Range("P7").Value = 144.9
$ws.Range("R7").Value = 145.1

# --- Row 8 ---
$ws.Range("N8").Value = 34.6
$ws.Range("P8").Value = 122.5
$ws.Range("R8").Value = 147.6

# --- Row 9 ---
$ws.Range("N9").Value = 143.55
$ws.Range("P9").Value = 186
$ws.Range("R9").Value = 233

# --- Row 23: fill in the previously blank Species_nu (B) column ---
$ws.Range("B23").Value = 5

# --- Row 39: recomputed regression statistics for Hypoatherina_temminckii (summary row) ---
$ws.Range("D39").Value = 0.019105
$ws.Range("E39").Value = 0.003273
$ws.Range("F39").Value = 2.686972
$ws.Range("G39").Value = 0.084889
$ws.Range("H39").Value = 0.3319
$ws.Range("I39").Value = 0.9672
$ws.Range("J39").Value = 0.9487459
$ws.Range("K39").Value = 1.047926
$ws.Range("L39").Value = "47 total"

# --- Row 43: recomputed regression statistics for Taeniamia_fucata (summary row) ---
$ws.Range("D43").Value = 0.0149764
$ws.Range("E43").Value = 0.0008158
$ws.Range("F43").Value = 3.184679

# --- Selection left on A53 after scrolling through the data ---
$ws.Activate()
$ws.Range("A53").Select()
